$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B2").Copy($ws.Range("A7:B7"))

$ws.Range("A7:B7").Select()
